# Insert one new data row at row 143 (pushes the existing rows 143-267 down
# to 144-268, growing the sheet's used range from A1:R267 to A1:R268).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(143).Insert()

$ws.Cells.Item(143, 1).Value  = 10
$ws.Cells.Item(143, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(143, 3).Value  = "La Araucanía"
$ws.Cells.Item(143, 4).Value  = 45040
$ws.Cells.Item(143, 5).Value  = 9
$ws.Cells.Item(143, 6).Value  = 100112012
$ws.Cells.Item(143, 7).Value  = "Espinaca"
$ws.Cells.Item(143, 8).Value  = "Sin especificar"
$ws.Cells.Item(143, 9).Value  = "Primera"
$ws.Cells.Item(143, 10).Value = 30
$ws.Cells.Item(143, 11).Value = 10000
$ws.Cells.Item(143, 12).Value = 10000
$ws.Cells.Item(143, 13).Value = 10000
$ws.Cells.Item(143, 14).Value = "$/docena de atados"
$ws.Cells.Item(143, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(143, 16).Value = 3333
$ws.Cells.Item(143, 17).Value = 3
$ws.Cells.Item(143, 18).Value = "Hortaliza"
